# edit.ps1
# Rewrites the "Quantum Realm" essay into "The Enigmatic Symphony of Politics..."
# per the supplied diff: title / byline / email swap, every body sentence
# rewritten from quantum-physics flavour text to a politics essay, a large
# new block of paragraphs appended to the intro, the summary rewritten, and
# a trailing empty paragraph added at the very end of the document.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

function Replace-Word($find, $replace) {
    # MatchWholeWord = $true so we don't clobber "computing"/"computation" etc.
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

$br = [string][char]11   # manual line break -> serializes as <w:br/>

# --- Title, byline, e-mail ---------------------------------------------------

Replace-Text "The Allure of the Quantum Realm" `
    "The Enigmatic Symphony of Politics: Navigating the Complex World of Power and Influence"

Replace-Text "Genevieve Smith" "Eleanor Winters"

Replace-Text "smithgenevieve@gmail" "eleanor"
Replace-Word "com" "winters@validedu"
Replace-Text "winters@validedu" "winters@validedu.org"

# --- Intro paragraph: three original sentence groups rewritten -------------

$p5_old_seg1 = "Journey with us into the captivating realm of quantum physics, where particles dance in a harmonious waltz of uncertainty and probability. This enigmatic realm, once confined to theoretical musings, is now seeping into our reality, promising awe-inspiring technologies that defy classical intuition. From the enigmatic world of quantum computing to the nascent field of quantum cryptography, we stand at the precipice of a paradigm shift, where the ethereal fabric of quantum mechanics is woven into the tapestry of our technological landscape."
$p5_new_seg1 = "In the realm of human affairs, there exists an intricate symphony of power, influence, and decision-making that we call politics. This vast and dynamic world of governance, leadership, and societal interactions shapes the very fabric of our societies, affecting every aspect of our lives, from the policies that govern us to the leaders who represent us. To navigate this complex landscape effectively, it is imperative that we understand the fundamental principles of politics, its historical evolution, and its profound impact on our daily lives."
Replace-Text $p5_old_seg1 $p5_new_seg1

$p5_old_seg2 = "Step into the arena of quantum computing, where information dances in the ethereal realm of quantum bits, also known as qubits. Unlike their classical counterparts, qubits wield the uncanny ability to exist in a superposition of states, pirouette-ing through a ballet of possibilities. This intoxicating dance grants quantum computers the potency to tackle conundrums that confound their classical brethren, pioneering solutions to intractable problems in cryptography, optimization, and simulation, unveiling secrets hidden within the labyrinthine pathways of computation."
$p5_new_seg2 = "Politics, at its core, is the art of resolving conflicts and allocating resources within a society. It involves the formulation and implementation of policies, the establishment of laws and regulations, and the distribution of power among various institutions and individuals. Through political processes, we collectively determine how we want to live together, what values we hold dear, and how we can create a just and equitable society for all."
Replace-Text $p5_old_seg2 $p5_new_seg2

$p5_old_seg3 = "Venture into the clandestine realm of quantum cryptography, where information cloaks itself in the enigmatic embrace of quantum mechanics, creating an impregnable shield against eavesdropping ears. This quantum cloak harnesses the inherent fragility of quantum information, orchestrating a symphony of particles that evokes alarm at the slightest touch of an unintended observer. With quantum cryptography as our guardian, we can forge unbreakable codes, ensuring the sanctity of our secrets in a world where data breaches are an incessant threat."
$p5_new_seg3 = "As we delve into the study of politics, we are confronted with a tapestry of historical events, political theories, and ideological debates that have shaped our current political landscape. From the ancient Greek city-states to the modern nation-states, from the rise and fall of empires to the emergence of global governance, politics has been an ever-evolving field, constantly adapting to changing circumstances and societal needs. By understanding the historical context of politics, we gain a deeper appreciation for the challenges and opportunities that lie ahead."
Replace-Text $p5_old_seg3 $p5_new_seg3

# --- New material appended to the end of that same paragraph ---------------

$p5Para = $d.Paragraphs.Item(5)

$newBlock = $br + $br + "Introduction Continued:" + $br + $br + `
    "Furthermore, the study of politics provides us with a framework for analyzing and understanding current events. By examining political institutions, policies, and decision-making processes, we can develop a critical perspective on the world around us. We can identify the various actors and interests at play, assess the potential impact of different policies, and engage in informed debates about the direction of our society. Politics is not merely an abstract concept; it is a living, breathing force that shapes our communities, our economies, and our planet." + `
    $br + $br + `
    "Politics is a multi-faceted subject that encompasses a wide range of topics, from the intricacies of international relations to the challenges of local governance. It involves the study of political systems, ideologies, public policy, and the role of citizens in a democracy. Through political engagement, we have the power to influence the decisions that affect our lives and to hold our leaders accountable"

$insertPoint = $d.Range($p5Para.Range.End - 1, $p5Para.Range.End - 1)
$insertPoint.InsertAfter($newBlock)
$insertedRange = $d.Range($p5Para.Range.End - 1 - $newBlock.Length, $p5Para.Range.End - 1)
$insertedRange.Font.Name = "Aptos"
$insertedRange.Font.Size = 12
$insertedRange.Font.Color = 0

# --- Summary paragraph -------------------------------------------------------

$p7_old = "The captivating realm of quantum physics dances on the boundary of our perception, blurring the line between theory and reality. Quantum computing, the nascent field of quantum cryptography, and the burgeoning world of quantum sensing are transforming our technological landscape. From decoding intricate problems to securing our digital realm, the quantum realm is reshaping our perception of what's possible. As we delve deeper into this enigmatic realm, we are unveiling a new chapter in the human quest for knowledge and dominion over the forces that govern our universe."
$p7_new = "In this essay, we have explored the enigmatic symphony of politics, a complex world of power, influence, and decision-making that profoundly impacts our lives. By understanding the fundamental principles of politics, its historical evolution, and its current manifestations, we gain a deeper appreciation for the challenges and opportunities that lie ahead. The study of politics equips us with the knowledge and skills necessary to navigate the intricacies of governance, to engage in informed debates, and to shape the future of our societies. As active citizens, we have the responsibility to participate in the political process, to hold our leaders accountable, and to strive for a more just and equitable world for all."
Replace-Text $p7_old $p7_new

# --- Trailing empty paragraph added at the very end of the document --------

$d.Content.InsertParagraphAfter() | Out-Null

Write-Output "edit.ps1 completed"
